# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Rebuilds the "Periodo Mora" detail table (rows 16-69) on Hoja1 so that the
# two workers (LIBIA MARTINEZ MEDINA / 33102376 and ROSIRIS TORRES TORRES /
# 45580363) are interleaved row-by-row across their full set of overdue
# periods (1810 .. 2102), instead of being listed as two separate blocks.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Hoja1")

# row, Tipo Doc, N Doc Trabajador, Nombre Trabajador, Periodo Mora, Valor Mora, Salario Basico
$data = @(
    @(16, "CC", "33102376", "LIBIA MARTINEZ MEDINA", "1810", 31249, 781242),
    @(17, "CC", "33102376", "LIBIA MARTINEZ MEDINA", "1811", 31249, 781242),
    @(18, "CC", "33102376", "LIBIA MARTINEZ MEDINA", "1812", 31249, 781242),
    @(19, "CC", "33102376", "LIBIA MARTINEZ MEDINA", "1901", 31249, 781242),
    @(20, "CC", "33102376", "LIBIA MARTINEZ MEDINA", "1902", 31249, 781242),
    @(21, "CC", "45580363", "ROSIRIS TORRES TORRES", "1902", 48000, 1200000),
    @(22, "CC", "33102376", "LIBIA MARTINEZ MEDINA", "1903", 31249, 781242),
    @(23, "CC", "45580363", "ROSIRIS TORRES TORRES", "1903", 48000, 1200000),
    @(24, "CC", "33102376", "LIBIA MARTINEZ MEDINA", "1904", 31249, 781242),
    @(25, "CC", "45580363", "ROSIRIS TORRES TORRES", "1904", 48000, 1200000),
    @(26, "CC", "33102376", "LIBIA MARTINEZ MEDINA", "1905", 31249, 781242),
    @(27, "CC", "45580363", "ROSIRIS TORRES TORRES", "1905", 48000, 1200000),
    @(28, "CC", "33102376", "LIBIA MARTINEZ MEDINA", "1906", 31249, 781242),
    @(29, "CC", "45580363", "ROSIRIS TORRES TORRES", "1906", 48000, 1200000),
    @(30, "CC", "33102376", "LIBIA MARTINEZ MEDINA", "1907", 31249, 781242),
    @(31, "CC", "45580363", "ROSIRIS TORRES TORRES", "1907", 48000, 1200000),
    @(32, "CC", "33102376", "LIBIA MARTINEZ MEDINA", "1908", 31249, 781242),
    @(33, "CC", "45580363", "ROSIRIS TORRES TORRES", "1908", 48000, 1200000),
    @(34, "CC", "33102376", "LIBIA MARTINEZ MEDINA", "1909", 31249, 781242),
    @(35, "CC", "45580363", "ROSIRIS TORRES TORRES", "1909", 48000, 1200000),
    @(36, "CC", "33102376", "LIBIA MARTINEZ MEDINA", "1910", 31249, 781242),
    @(37, "CC", "45580363", "ROSIRIS TORRES TORRES", "1910", 48000, 1200000),
    @(38, "CC", "33102376", "LIBIA MARTINEZ MEDINA", "1911", 31249, 781242),
    @(39, "CC", "45580363", "ROSIRIS TORRES TORRES", "1911", 48000, 1200000),
    @(40, "CC", "33102376", "LIBIA MARTINEZ MEDINA", "1912", 31249, 781242),
    @(41, "CC", "45580363", "ROSIRIS TORRES TORRES", "1912", 48000, 1200000),
    @(42, "CC", "33102376", "LIBIA MARTINEZ MEDINA", "2001", 31249, 781242),
    @(43, "CC", "45580363", "ROSIRIS TORRES TORRES", "2001", 48000, 1200000),
    @(44, "CC", "33102376", "LIBIA MARTINEZ MEDINA", "2002", 31249, 781242),
    @(45, "CC", "45580363", "ROSIRIS TORRES TORRES", "2002", 48000, 1200000),
    @(46, "CC", "33102376", "LIBIA MARTINEZ MEDINA", "2003", 31249, 781242),
    @(47, "CC", "45580363", "ROSIRIS TORRES TORRES", "2003", 48000, 1200000),
    @(48, "CC", "33102376", "LIBIA MARTINEZ MEDINA", "2004", 31249, 781242),
    @(49, "CC", "45580363", "ROSIRIS TORRES TORRES", "2004", 48000, 1200000),
    @(50, "CC", "33102376", "LIBIA MARTINEZ MEDINA", "2005", 31249, 781242),
    @(51, "CC", "45580363", "ROSIRIS TORRES TORRES", "2005", 48000, 1200000),
    @(52, "CC", "33102376", "LIBIA MARTINEZ MEDINA", "2006", 31249, 781242),
    @(53, "CC", "45580363", "ROSIRIS TORRES TORRES", "2006", 48000, 1200000),
    @(54, "CC", "33102376", "LIBIA MARTINEZ MEDINA", "2007", 31249, 781242),
    @(55, "CC", "45580363", "ROSIRIS TORRES TORRES", "2007", 48000, 1200000),
    @(56, "CC", "33102376", "LIBIA MARTINEZ MEDINA", "2008", 31249, 781242),
    @(57, "CC", "45580363", "ROSIRIS TORRES TORRES", "2008", 48000, 1200000),
    @(58, "CC", "33102376", "LIBIA MARTINEZ MEDINA", "2009", 31249, 781242),
    @(59, "CC", "45580363", "ROSIRIS TORRES TORRES", "2009", 48000, 1200000),
    @(60, "CC", "33102376", "LIBIA MARTINEZ MEDINA", "2010", 31249, 781242),
    @(61, "CC", "45580363", "ROSIRIS TORRES TORRES", "2010", 48000, 1200000),
    @(62, "CC", "33102376", "LIBIA MARTINEZ MEDINA", "2011", 31249, 781242),
    @(63, "CC", "45580363", "ROSIRIS TORRES TORRES", "2011", 48000, 1200000),
    @(64, "CC", "33102376", "LIBIA MARTINEZ MEDINA", "2012", 31249, 781242),
    @(65, "CC", "45580363", "ROSIRIS TORRES TORRES", "2012", 48000, 1200000),
    @(66, "CC", "33102376", "LIBIA MARTINEZ MEDINA", "2101", 31249, 781242),
    @(67, "CC", "45580363", "ROSIRIS TORRES TORRES", "2101", 48000, 1200000),
    @(68, "CC", "33102376", "LIBIA MARTINEZ MEDINA", "2102", 22916, 781242),
    @(69, "CC", "45580363", "ROSIRIS TORRES TORRES", "2102", 35200, 1200000)
)

foreach ($r in $data) {
    $row = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]        # B Tipo Doc Trabajador
    $ws.Cells.Item($row, 3).Value = $r[2]        # C N Doc Trabajador
    $ws.Cells.Item($row, 4).Value = $r[3]        # D Nombre Trabajador
    $ws.Cells.Item($row, 5).Value = $r[4]        # E Periodo Mora
    $ws.Cells.Item($row, 6).Value = $r[5]        # F Valor Mora
    $ws.Cells.Item($row, 7).Value = $r[6]        # G Salario Basico
}
